# [Feat] REMOVE 패킷 추가
# - 접속 해제 시 캐릭터를 지우도록 REMOVE 패킷을 추가함
# - remove 메시지 전송 시 방 id가 아닌 플레이어 id로 방 정보를 찾던 버그 수정
#
# Inserts a new "SC_REMOVE_CHARACTER" message definition (playerID / uint32)
# above the existing "PlayerInfo" block, keeping the blank separator-row
# layout used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")

# Push the existing "PlayerInfo" block (and everything below it) down by two
# rows: one row for the new message row, one row to keep it as its own
# blank-separated block, same as every other message group on this sheet.
$ws.Rows("44:45").Insert()

# New message: SC_REMOVE_CHARACTER(playerID: uint32)
$ws.Range("A44").Value = "SC_REMOVE_CHARACTER"
$ws.Range("B44").Value = "playerID"
$ws.Range("C44").Value = "uint32"
$ws.Range("D44").Value = "플레이어의 id"

# Restore view: scroll down a bit and select A45 (matches author's saved view).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A45").Select()
